# Rename the worksheet from "Sheet1" to "Combined"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Combined"

# Remove all cell comments on the sheet (this also drops the now-unused
# VML legacy-drawing part/relationship that Excel keeps solely to render
# the comment boxes).
foreach ($sheet in $wb.Worksheets) {
    while ($sheet.Comments.Count -gt 0) {
        $sheet.Comments.Item(1).Delete()
    }
}
